$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty values for row 9
$ws.Range("B9").Value = 10852.325999999999
$ws.Range("C9").Value = 10000
$ws.Range("D9").Value = 659
$ws.Range("E9").Value = 178790202

# Update the active selection shown in the saved view
$ws.Range("C17").Select()
